$wb = $excel.ActiveWorkbook

# --- Update status text: "Ready for handoff" -> "In Translation" ---
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        if ("Ready for handoff" -eq $cell.Value()) {
            $cell.Value = "In Translation"
        }
    }
}

# --- Update column widths ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E1").ColumnWidth = 13.4101845877511
$overview.Range("F1").ColumnWidth = 13.4101845877511

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C1").ColumnWidth = 13.4101845877511

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C1").ColumnWidth = 13.4101845877511
